$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets: Sheet1 -> Daily, Sheet2 -> Fund
$ws1.Name = "Daily"
$ws2.Name = "Fund"

# --- Fund sheet (was Sheet2): fill in the Funds Input table ---

# Header row (row 2): Date | Debit | Credit | (blank) | Total
$ws2.Range("B2").Value = "Date"
$ws2.Range("D2").Value = "Debit"
$ws2.Range("F2").Value = "Credit"
$ws2.Range("H2").Value = "Total"

$ws2.Range("B2").Font.Bold = $true
$ws2.Range("D2").Font.Bold = $true
$ws2.Range("F2").Font.Bold = $true
$ws2.Range("G2").Font.Bold = $true
$ws2.Range("H2").Font.Bold = $true

$ws2.Rows.Item(2).HorizontalAlignment = -4108
$ws2.Range("A2").Clear()
$ws2.Range("C2").Clear()
$ws2.Range("E2").Clear()

# Data rows
$ws2.Range("B4").Value = 43019
$ws2.Range("B4").NumberFormat = "d-mmm"
$ws2.Range("D4").Value = 101
$ws2.Range("H4").Formula = "=D4"

$ws2.Range("B5").Value = 43019
$ws2.Range("B5").NumberFormat = "d-mmm"
$ws2.Range("D5").Value = 50000
$ws2.Range("H5").Formula = "=H4+D5-F5"

$ws2.Range("B6").Value = 43020
$ws2.Range("B6").NumberFormat = "d-mmm"
$ws2.Range("D6").Value = 50000

$ws2.Range("B7").Value = 43040
$ws2.Range("B7").NumberFormat = "d-mmm"
$ws2.Range("D7").Value = 20000

$ws2.Range("B8").Value = 43045
$ws2.Range("B8").NumberFormat = "d-mmm"
$ws2.Range("D8").Value = 75000

$ws2.Range("B9").Value = 43052
$ws2.Range("B9").NumberFormat = "d-mmm"
$ws2.Range("D9").Value = 35000

# Running-total formula, entered once across the block so it is stored as a
# single shared formula (H6:H9), same as Excel does when you fill a formula
# down a contiguous range.
$ws2.Range("H6:H9").Formula = "=H5+D6-F6"

# Column widths tweaked on Fund sheet
$ws2.Columns.Item(4).ColumnWidth = 10.583333333333334
$ws2.Columns.Item(6).ColumnWidth = 10.416666666666666
$ws2.Columns.Item(8).ColumnWidth = 10.333333333333334

# --- Selection / active-sheet bookkeeping ---
# Daily sheet's selection moves from G16 to A4
$ws1.Range("A4").Select() | Out-Null
# Fund sheet becomes the active tab, with H11 selected
$ws2.Range("H11").Select() | Out-Null
